$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Worksheet")

# Fill in Sprint 4 (row 10) time values for each contributor column
$ws.Range("B10").Value = 0.35694444444444445
$ws.Range("C10").Value = 0.50069444444444444
$ws.Range("D10").Value = 0.48541666666666666
$ws.Range("E10").Value = 0.34791666666666665

# New number format style applied to B10/D10/E10 (numFmtId 20 "h:mm"),
# keep existing font/alignment formatting intact
$ws.Range("B10").NumberFormat = "h:mm"
$ws.Range("D10").NumberFormat = "h:mm"
$ws.Range("E10").NumberFormat = "h:mm"

# Update sheet view: scroll the window so row 4 is at the top, then
# select the cell that ends up active/highlighted
$ws.Activate()
$excel.Goto($ws.Range("A4"), $true)
$ws.Range("C19").Select()
